# Update "想去人数" (column F) counts on both the "展览" and "全部类型"
# sheets to reflect freshly-scraped totals (gh-pages data regeneration).

$wb = $excel.ActiveWorkbook

# Row -> new F value mapping (same updates apply identically to both sheets).
$updates = @{
    3  = 2198
    5  = 13058
    6  = 73
    11 = 977
    12 = 13754
    13 = 14317
    25 = 5380
    29 = 15
    30 = 25
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
